# "More data for optimized helicopter" -- append a second data block
# (rows 13-17) below the existing experiment table, labelled with two
# new notes: "After Optimizing" and "same as above".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: section label
$ws.Range("A13").Value = "After Optimizing"

# Row 14: full data row (width, length, nose, time1-3, dist1-3)
$ws.Range("A14").Value = 2.875
$ws.Range("B14").Value = 4.75
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 2.03
$ws.Range("E14").Value = 1.86
$ws.Range("F14").Value = 2.13
$ws.Range("I14").Value = 15.5
$ws.Range("J14").Value = 5.25
$ws.Range("K14").Value = 13

# Row 15: label row, partial data (D,E time + I,J dist only)
$ws.Range("A15").Value = "same as above"
$ws.Range("D15").Value = 2.05
$ws.Range("E15").Value = 2.1
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 5

# Row 16: full data row
$ws.Range("A16").Value = 5.5
$ws.Range("B16").Value = 6
$ws.Range("C16").Value = 0.5
$ws.Range("D16").Value = 1.8
$ws.Range("E16").Value = 2.11
$ws.Range("F16").Value = 1.83
$ws.Range("I16").Value = 27.5
$ws.Range("J16").Value = 9
$ws.Range("K16").Value = 21

# Row 17: trailing label only
$ws.Range("A17").Value = "same as above"

# Leave the selection where the author's cursor ended up after entry
$ws.Range("G19").Select()
